$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Volume 30  Number NN" text (A8), edit just the number run ---
$volChars = $ws.Range("A8").Characters(21, 2)
$volChars.Text = "48"

# --- Update "Report Covering the Week  MM/DD/YYYY  Through  MM/DD/YYYY" (C9) ---
$weekCell = $ws.Range("C9")
$weekCell.Characters(27, 10).Text = "11/27/2023"
$weekCell.Characters(48, 10).Text = "12/3/2023"

# --- Update weekly crime-stat data cells (rows 14-21, 24-27) ---
$ws.Range("N14").Value = -50
$ws.Range("D15").Value = 1
$ws.Range("D15").NumberFormat = '#,##0'
$ws.Range("E15").Value = -100
$ws.Range("E15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("G15").Value = 2
$ws.Range("J15").Value = 4
$ws.Range("K15").Value = 100
$ws.Range("N15").Value = -38.461538461538
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 2
$ws.Range("F16").Value = 3
$ws.Range("G16").Value = 6
$ws.Range("H16").Value = -50
$ws.Range("I16").Value = 56
$ws.Range("J16").Value = 79
$ws.Range("K16").Value = -29.113924050632
$ws.Range("L16").Value = 107.407407407407
$ws.Range("M16").Value = -28.205128205128
$ws.Range("N16").Value = -81.758957654723
$ws.Range("C17").Value = 6
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 10
$ws.Range("G17").Value = 9
$ws.Range("H17").Value = 11.111111111111
$ws.Range("I17").Value = 95
$ws.Range("J17").Value = 95
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 30.136986301369
$ws.Range("M17").Value = 93.877551020408
$ws.Range("N17").Value = -24.603174603174
$ws.Range("C18").Value = 11
$ws.Range("D18").Value = 21
$ws.Range("E18").Value = -47.619047619047
$ws.Range("F18").Value = 38
$ws.Range("G18").Value = 50
$ws.Range("H18").Value = -24
$ws.Range("I18").Value = 302
$ws.Range("J18").Value = 279
$ws.Range("K18").Value = 8.243727598566
$ws.Range("L18").Value = 27.966101694915
$ws.Range("M18").Value = 28.510638297872
$ws.Range("N18").Value = -67.102396514161
$ws.Range("C19").Value = 14
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = 40
$ws.Range("F19").Value = 41
$ws.Range("G19").Value = 49
$ws.Range("H19").Value = -16.326530612244
$ws.Range("I19").Value = 585
$ws.Range("J19").Value = 612
$ws.Range("K19").Value = -4.411764705882
$ws.Range("L19").Value = 65.722379603399
$ws.Range("M19").Value = 69.565217391304
$ws.Range("N19").Value = 10.377358490566
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 7
$ws.Range("E20").Value = -42.857142857142
$ws.Range("F20").Value = 19
$ws.Range("G20").Value = 14
$ws.Range("H20").Value = 35.714285714285
$ws.Range("I20").Value = 179
$ws.Range("J20").Value = 112
$ws.Range("K20").Value = 59.821428571428
$ws.Range("L20").Value = 179.6875
$ws.Range("M20").Value = 39.84375
$ws.Range("N20").Value = -94.165580182529
$ws.Range("C21").Value = 37
$ws.Range("D21").Value = 44
$ws.Range("E21").Value = -15.909090909090
$ws.Range("F21").Value = 111
$ws.Range("G21").Value = 130
$ws.Range("H21").Value = -14.615384615384
$ws.Range("I21").Value = 1227
$ws.Range("J21").Value = 1181
$ws.Range("K21").Value = 3.895004233700
$ws.Range("L21").Value = 61.873350923482
$ws.Range("M21").Value = 45.89774078478
$ws.Range("N21").Value = -75.291985501409
$ws.Range("C24").Value = 7
$ws.Range("D24").Value = 11
$ws.Range("E24").Value = -36.363636363636
$ws.Range("F24").Value = 47
$ws.Range("G24").Value = 52
$ws.Range("H24").Value = -9.615384615384
$ws.Range("I24").Value = 537
$ws.Range("J24").Value = 711
$ws.Range("K24").Value = -24.472573839662
$ws.Range("L24").Value = 2.091254752851
$ws.Range("M24").Value = 28.162291169451
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 0
$ws.Range("G25").Value = 16
$ws.Range("H25").Value = 31.25
$ws.Range("I25").Value = 214
$ws.Range("J25").Value = 227
$ws.Range("K25").Value = -5.726872246696
$ws.Range("L25").Value = 47.586206896551
$ws.Range("M25").Value = 32.098765432098
$ws.Range("D26").Value = 1
$ws.Range("D26").NumberFormat = '#,##0'
$ws.Range("E26").Value = -100
$ws.Range("E26").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("G26").Value = 2
$ws.Range("J26").Value = 5
$ws.Range("K26").Value = 80
$ws.Range("C27").Value = 1
$ws.Range("C27").NumberFormat = '#,##0'
$ws.Range("F27").Value = 2
$ws.Range("I27").Value = 18
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = -18.181818181818
